$d = $word.ActiveDocument

# 1. Title / heading text (appears twice: H1 heading and bold summary line)
$d.Content.Find.Execute("Play Eye of Horus Megaways Free: A Modernized Ancient Egypt Slot", $true, $false, $false, $false, $false, $true, 1, $false, "Play Eye of Horus Megaways for Free - Review", 2)

# 2. "What we like" bullet list
$d.Content.Find.Execute("Up to 15,625 ways to win", $true, $false, $false, $false, $false, $true, 1, $false, "Megaways gaming system adds excitement to gameplay", 2)
$d.Content.Find.Execute("Potential jackpot of 10,000 times the bet", $true, $false, $false, $false, $false, $true, 1, $false, "Impressive winning potential with up to 10,000x bet", 2)
$d.Content.Find.Execute("Wild symbol behaves as an expandable symbol during free spins", $true, $false, $false, $false, $false, $true, 1, $false, "High-quality graphics and modernized design", 2)
$d.Content.Find.Execute("Thematic symbols related to Ancient Egypt culture", $true, $false, $false, $false, $false, $true, 1, $false, "Free spins feature with expandable Wild symbol for big wins", 2)

# 3. "What we don't like" bullet list
$d.Content.Find.Execute("No bonus games outside of free spins", $true, $false, $false, $false, $false, $true, 1, $false, "Limited extra features beyond the free spins", 2)
$d.Content.Find.Execute("Non-themed regular symbols have low payout", $true, $false, $false, $false, $false, $true, 1, $false, "Theme of Ancient Egypt may feel overused for some players", 2)

# 4. Meta description (italic) text
$d.Content.Find.Execute("Read our review of Eye of Horus Megaways and play it for free. Experience the modern version of an Ancient Egypt slot with up to 15,625 ways to win.", $true, $false, $false, $false, $false, $true, 1, $false, "Find out all about Eye of Horus Megaways and play for free. Review of gameplay, features, and winning potential.", 2)
